$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$caseFilesQuery = @'
SELECT 
    DISTINCT cf.file_name AS "File Name",
    CASE
        WHEN cf.file_name LIKE '%.bai' THEN 'bai'
        WHEN cf.file_name LIKE '%.bam' THEN 'bam'
        WHEN cf.file_name LIKE '%.csv' THEN 'csv'
        WHEN cf.file_name LIKE '%.doc' THEN 'doc'
        WHEN cf.file_name LIKE '%.docx' THEN 'docx'
         WHEN cf.file_name LIKE '%.fastq.gz' THEN 'fastq'
        WHEN cf.file_name LIKE '%.vcf.gz' THEN 'vcf'
        WHEN cf.file_name LIKE '%.tsv.gz' THEN 'tsv'
        WHEN cf.file_name LIKE '%.pdf' THEN 'pdf'
        WHEN cf.file_name LIKE '%.rtf' THEN 'rtf'
        WHEN cf.file_name LIKE '%.tbi' THEN 'tbi'
        WHEN cf.file_name LIKE '%.tif' THEN 'tif'
        WHEN cf.file_name LIKE '%.xls' THEN 'xls'
        WHEN cf.file_name LIKE '%.xlsx' THEN 'xlsx'
        ELSE 'Unknown'
    END AS "Format",
    cf.file_type AS "File Type",
    CASE     
    WHEN cf.file_size >= 1024 * 1024 * 1024 THEN 
        CASE 
            WHEN ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 2) = CAST(ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT) 
            THEN CAST(CAST(ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT) AS TEXT) || ' GB'
            ELSE ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 2) || ' GB'
        END
    WHEN cf.file_size >= 1024 * 1024 THEN 
        CASE 
            WHEN ROUND(cf.file_size / (1024.0 * 1024.0), 2) = CAST(ROUND(cf.file_size / (1024.0 * 1024.0), 0) AS INT) 
            THEN CAST(CAST(ROUND(cf.file_size / (1024.0 * 1024.0), 0) AS INT) AS TEXT) || ' MB'
            ELSE ROUND(cf.file_size / (1024.0 * 1024.0), 2) || ' MB'
        END
    WHEN cf.file_size >= 1024 THEN 
        CASE 
            WHEN ROUND(cf.file_size / 1024.0, 2) = CAST(ROUND(cf.file_size / 1024.0, 0) AS INT) 
            THEN CAST(CAST(ROUND(cf.file_size / 1024.0, 0) AS INT) AS TEXT) || ' KB'
            ELSE ROUND(cf.file_size / 1024.0, 2) || ' KB'
        END
    ELSE 
        CASE 
            WHEN ROUND(cf.file_size, 2) = CAST(ROUND(cf.file_size, 0) AS INT) 
            THEN CAST(CAST(ROUND(cf.file_size, 0) AS INT) AS TEXT) || ' Bytes'
            ELSE ROUND(cf.file_size, 2) || ' Bytes'
        END
END AS "Size",
    'sample' AS "Association",
    cf.file_description AS "Description",
    smp.sample_id AS "Sample ID",
    c.case_record_id AS "Case ID",
    dmg.breed AS "Breed",
    diag.disease_term AS "Diagnosis"
FROM 
    df_case_file cf
JOIN 
    df_sample smp ON cf."sample.sample_id" = smp.sample_id
JOIN 
    df_case c ON smp."case.case_record_id" = c.case_record_id
JOIN 
    df_study st ON c."study.clinical_study_designation" = st.clinical_study_designation
JOIN 
    df_program p ON st."program.program_acronym" = p.program_acronym
JOIN 
    df_demographic dmg ON dmg."case.case_record_id" = c.case_record_id
JOIN 
    df_diagnosis diag ON diag."case.case_record_id" = c.case_record_id
JOIN 
    df_enrollment enr ON enr."case.case_record_id" = c.case_record_id
JOIN 
    df_publication pub ON pub."study.clinical_study_designation" = st.clinical_study_designation
LEFT JOIN 
    df_study_file sf ON sf."study.clinical_study_designation" = st.clinical_study_designation
WHERE
    st.clinical_study_designation = 'COTC022' AND smp.summarized_sample_type = 'Primary Malignant Tumor Tissue'
ORDER BY 
    cf.file_name ASC
LIMIT 100;
'@

$ws.Range("B4").Value = $caseFilesQuery

# The row was already wrapped at Excel's maximum row height before the edit;
# keep it pinned there instead of letting autofit grow it further.
$ws.Rows.Item(4).RowHeight = 409.6

# Move the active selection the same way the author's Excel session left it
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
